$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" summary sheet: insert a new row for 2022-Q4 above the existing
#    2021-Q2 row, pushing 2021-Q2 down to row 3.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Push the existing 2021-Q2 summary row (row 2) down to row 3, copying A2's
# style onto A3 first so the new row keeps the same look.
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2021-Q2"
$wsTotal.Range("C3").Value = 4
$wsTotal.Range("D3").Value = 0.25

# Overwrite row 2 with the new 2022-Q4 summary entry (A2 keeps its value/style).
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 3
$wsTotal.Range("D2").Value = 0.27

# ---------------------------------------------------------------------------
# 2) Duplicate the existing "2021-Q2" detail sheet so the original data is
#    preserved intact on its own tab, then repurpose the original tab (which
#    keeps its original sheetId/position) to hold the new 2022-Q4 detail
#    data. This reproduces: 总计, 2022-Q4 (new), 2021-Q2 (old, moved).
# ---------------------------------------------------------------------------
$sheetOld = $wb.Worksheets.Item("2021-Q2")
$sheetOld.Copy($null, $sheetOld)
$sheetDup = $wb.Worksheets.Item("2021-Q2 (2)")

$sheetOld.Name = "2022-Q4"
$sheetDup.Name = "2021-Q2"

# ---------------------------------------------------------------------------
# 3) Populate the (now empty-able) "2022-Q4" sheet with the new fund detail
#    data, replacing whatever it held before.
# ---------------------------------------------------------------------------
$ws2022 = $wb.Worksheets.Item("2022-Q4")
$ws2022.Cells.Clear()

$ws2022.Range("B1").Value = "基金代码"
$ws2022.Range("C1").Value = "基金名称"
$ws2022.Range("D1").Value = "基金规模"
$ws2022.Range("E1").Value = "股票总仓位"
$ws2022.Range("F1").Value = "仓位占比"
$ws2022.Range("G1").Value = "持有市值(亿元)"
$ws2022.Range("H1").Value = "仓位排名"

$ws2022.Range("A2").Value = 0
$ws2022.Range("B2").Value = "'000593"
$ws2022.Range("C2").Value = "易方达标普全球高端消费品指数增强（QDII）美元现汇"
$ws2022.Range("D2").Value = "'2.30"
$ws2022.Range("E2").Value = "'93.71"
$ws2022.Range("F2").Value = "'3.88"
$ws2022.Range("G2").Value = "'0.0892"
$ws2022.Range("H2").Value = 9

$ws2022.Range("A3").Value = 1
$ws2022.Range("B3").Value = "'005676"
$ws2022.Range("C3").Value = "易方达标普全球高端消费品指数增强C（QDII）人民币"
$ws2022.Range("D3").Value = "'2.30"
$ws2022.Range("E3").Value = "'93.71"
$ws2022.Range("F3").Value = "'3.88"
$ws2022.Range("G3").Value = "'0.0892"
$ws2022.Range("H3").Value = 9

$ws2022.Range("A4").Value = 2
$ws2022.Range("B4").Value = "'118002"
$ws2022.Range("C4").Value = "易方达标普全球高端消费品指数增强A（QDII）人民币"
$ws2022.Range("D4").Value = "'2.30"
$ws2022.Range("E4").Value = "'93.71"
$ws2022.Range("F4").Value = "'3.88"
$ws2022.Range("G4").Value = "'0.0892"
$ws2022.Range("H4").Value = 9

# Re-apply the header / index-column style (matches the bold centred style
# already used for the "总计" sheet's headers) since Cells.Clear() wiped it.
$wsTotal.Range("B1").Copy()
$ws2022.Range("B1:H1").PasteSpecial(-4122)

$wsTotal.Range("A2").Copy()
$ws2022.Range("A2:A4").PasteSpecial(-4122)
$ws2022.Range("A2").Value = 0
$ws2022.Range("A3").Value = 1
$ws2022.Range("A4").Value = 2

$wsTotal.Range("A1").Select()
